$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old hyperlink & cell at D512 (moves to D539) ---
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 512) {
        $h.Delete()
    }
}
$ws.Range("D512").Clear()

# --- Populate new rows 539-597 (column order A,B,C,D matches original translator entry order so shared-string indices line up) ---
$ws.Range("A539").Value = 143
$ws.Range("B539").Value = 'Volume Settings'
$ws.Range("C539").Value = 'Impostazioni del volume'

# --- Add hyperlink at D539 (same target as before, moved down from D512) ---
$ws.Hyperlinks.Add($ws.Range("D539"), "https://www.deepl.com/translator", "", "", "https://www.deepl.com/translator")
$ws.Range("D539").Value = "DeepL Translate: The world's most accurate translator"
$ws.Range("D539").Style = $ws.Range("F1").Style

$ws.Range("B540").Value = 'Device name:'
$ws.Range("C540").Value = 'Nome del dispositivo:'
$ws.Range("B541").Value = 'Volume label:'
$ws.Range("C541").Value = 'Etichetta del volume:'
$ws.Range("B542").Value = 'Path:'
$ws.Range("C542").Value = 'Percorso:'
$ws.Range("B543").Value = 'Read/write'
$ws.Range("C543").Value = 'Lettura/scrittura'
$ws.Range("B544").Value = 'Bootable'
$ws.Range("C544").Value = 'Avviabile'
$ws.Range("B545").Value = 'Boot priority:'
$ws.Range("C545").Value = 'Priorità di avvio:'
$ws.Range("B546").Value = 'Select Directory'
$ws.Range("C546").Value = 'Selezionare la directory'
$ws.Range("B547").Value = 'Select Archive or Plain File'
$ws.Range("C547").Value = 'Selezionare Archivio o File semplice'
$ws.Range("B548").Value = 'OK'
$ws.Range("C548").Value = 'OK'
$ws.Range("B549").Value = 'Cancel'
$ws.Range("C549").Value = 'Annullamento'
$ws.Range("B550").Value = 'Eject'
$ws.Range("C550").Value = 'Espulsione'
$ws.Range("A551").Value = 144
$ws.Range("B551").Value = 'When Active'
$ws.Range("C551").Value = 'Quando è attivo'
$ws.Range("B552").Value = 'Run at priority:'
$ws.Range("C552").Value = 'Eseguire con priorità:'
$ws.Range("B553").Value = 'Mouse uncaptured:'
$ws.Range("C553").Value = 'Topo non catturato:'
$ws.Range("B554").Value = 'Pause emulation'
$ws.Range("C554").Value = 'Pausa emulazione'
$ws.Range("B555").Value = 'Disable sound'
$ws.Range("C555").Value = 'Disattivare il suono'
$ws.Range("B556").Value = 'When Inactive'
$ws.Range("C556").Value = 'Quando è inattivo'
$ws.Range("B557").Value = 'Disable game controllers'
$ws.Range("C557").Value = 'Disattivare i controller di gioco'
$ws.Range("B558").Value = 'When Minimized'
$ws.Range("C558").Value = 'Quando è ridotto al minimo'
$ws.Range("B559").Value = 'File Extension Associations'
$ws.Range("C559").Value = 'Associazioni di estensioni di file'
$ws.Range("B560").Value = 'Associate all'
$ws.Range("C560").Value = 'Associare tutti'
$ws.Range("B561").Value = 'Deassociate all'
$ws.Range("C561").Value = 'Dissociare tutti i'
$ws.Range("A562").Value = 152
$ws.Range("B562").Value = 'WinUAE Debugger'
$ws.Range("C562").Value = 'Debugger WinUAE'
$ws.Range("A563").Value = 153
$ws.Range("B563").Value = 'Additional Information Settings'
$ws.Range("C563").Value = 'Informazioni aggiuntive Impostazioni'
$ws.Range("B564").Value = 'Path:'
$ws.Range("C564").Value = 'Percorso:'
$ws.Range("B565").Value = 'Link:'
$ws.Range("C565").Value = 'Link:'
$ws.Range("B566").Value = 'Category:'
$ws.Range("C566").Value = 'Categoria:'
$ws.Range("B567").Value = 'Tags:'
$ws.Range("C567").Value = 'Tag:'
$ws.Range("B568").Value = 'Ignore link'
$ws.Range("C568").Value = 'Ignorare il link'
$ws.Range("B569").Value = 'Autoload'
$ws.Range("C569").Value = 'Carica automatica'
$ws.Range("B570").Value = 'OK'
$ws.Range("C570").Value = 'OK'
$ws.Range("B571").Value = 'Cancel'
$ws.Range("C571").Value = 'Annullamento'
$ws.Range("A572").Value = 154
$ws.Range("B572").Value = 'Chipset'
$ws.Range("B573").Value = 'OCS [] Original chipset. A1000 and most A500s.'
$ws.Range("B574").Value = 'ECS Agnus [] Enhanced chipset (ECS Agnus chip only). CDTV and later A500 and A2000 hardware revisions.'
$ws.Range("B575").Value = 'Full ECS [] Full ECS chipset (ECS Agnus and ECS Denise chips). A500+, A600 and A3000.'
$ws.Range("B576").Value = 'AGA [] Advanced Graphics Architecture chipset. A1200, A4000 and CD32.'
$ws.Range("B577").Value = 'ECS Denise [] Enhanced chipset (ECS Denise chip only). Normally paired with ECS Agnus.'
$ws.Range("B578").Value = 'NTSC [] North American and Japanese display standard, 60Hz refresh rate. Other countries use PAL (50Hz. display refresh rate)'
$ws.Range("B579").Value = 'Cycle-exact (Full) [] The most compatible A500/A1200 emulation mode.'
$ws.Range("B580").Value = 'Cycle-exact (DMA/Memory accesses)'
$ws.Range("B581").Value = 'Chipset Extra:'
$ws.Range("B582").Value = 'Options'
$ws.Range("B583").Value = 'Keyboard connected'
$ws.Range("B584").Value = 'Subpixel display emulation'
$ws.Range("B585").Value = 'Immediate Blitter [] Faster but less compatible blitter emulation.'
$ws.Range("B586").Value = 'Wait for Blitter [] Compatibility hack for programs that don''t wait for the blitter correctly, causing graphics corruption if CPU is too fast.'
$ws.Range("B587").Value = 'Video port display hardware:'
$ws.Range("B588").Value = 'Monitor:'
$ws.Range("B589").Value = 'Collision Level'
$ws.Range("B590").Value = 'None [] Collision hardware emulation disabled.'
$ws.Range("B591").Value = 'Sprites only [] Emulate only sprite vs. sprite collisions.'
$ws.Range("B592").Value = 'Sprites and Sprites vs. Playfield [] Recommended collision emulation level.'
$ws.Range("B593").Value = 'Full [] 100% collision hardware emulation. Only very few games need this option. Slowest.'
$ws.Range("B594").Value = 'Genlock'
$ws.Range("B595").Value = 'Genlock connected [] Allow boot sequence to detect genlock.'
$ws.Range("B596").Value = 'Include alpha channel in screenshots and video captures.'
$ws.Range("B597").Value = 'Keep aspect ratio'

# --- Update selection / view ---
$ws.Range("B598").Select()
